$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: row index, C, D, E, F, G (new values per diff)
$rows = @(
    @(2, 1.511179703834972, 1.615307186802609, 1977, 0.3566017316017316, 31.38095238095238),
    @(3, 1.511971737413812, 1.616002815808326, 1975, 0.3562409812409812, 31.34920634920635),
    @(4, 1.511971737413812, 1.616002815808326, 1975, 0.3562409812409812, 31.34920634920635),
    @(5, 1.512444434566423, 1.616739276620547, 1973, 0.3558802308802309, 31.31746031746032),
    @(6, 1.512672818358493, 1.616709701541518, 1972, 0.3556998556998557, 31.3015873015873),
    @(7, 1.512672818358493, 1.616709701541518, 1972, 0.3556998556998557, 31.3015873015873),
    @(8, 1.512672818358493, 1.616709701541518, 1972, 0.3556998556998557, 31.3015873015873),
    @(9, 1.512672818358493, 1.616709701541518, 1972, 0.3556998556998557, 31.3015873015873),
    @(10, 1.512672818358493, 1.616709701541518, 1972, 0.3556998556998557, 31.3015873015873),
    @(11, 1.513929005258451, 1.617119499269849, 1971, 0.3555194805194805, 31.28571428571428),
    @(12, 1.513929005258451, 1.617119499269849, 1971, 0.3555194805194805, 31.28571428571428),
    @(13, 1.514310337456992, 1.617524494725992, 1970, 0.3553391053391053, 31.26984126984127),
    @(14, 1.516064366965211, 1.618889691572275, 1966, 0.3546176046176046, 31.20634920634921),
    @(15, 1.51547201939842, 1.617602851623264, 1964, 0.3542568542568543, 31.17460317460317),
    @(16, 1.513919357203425, 1.616075986497502, 1961, 0.3537157287157287, 31.12698412698413),
    @(17, 1.514874815463211, 1.617468869798642, 1957, 0.352994227994228, 31.06349206349206),
    @(18, 1.515748436013758, 1.618060899210338, 1955, 0.3526334776334776, 31.03174603174603),
    @(19, 1.510429358220765, 1.614010689049391, 1949, 0.3515512265512266, 30.93650793650794),
    @(20, 1.513172588032851, 1.614742969456763, 1946, 0.351010101010101, 30.88888888888889),
    @(21, 1.514664898520284, 1.615585039048565, 1943, 0.3504689754689755, 30.84126984126984),
    @(22, 1.514542169433464, 1.614784961076318, 1942, 0.3502886002886003, 30.82539682539683),
    @(23, 1.518173956024117, 1.615473723862625, 1938, 0.3495670995670996, 30.76190476190476),
    @(24, 1.512999259120564, 1.610220451286226, 1934, 0.3488455988455988, 30.6984126984127),
    @(25, 1.509901782911796, 1.60362754135532, 1925, 0.3472222222222222, 30.55555555555556),
    @(26, 1.498443291292488, 1.599852147852569, 1918, 0.345959595959596, 30.44444444444444),
    @(27, 1.500660291190787, 1.599873007921771, 1908, 0.3441558441558442, 30.28571428571428),
    @(28, 1.505725214943513, 1.600770805704966, 1903, 0.3432539682539683, 30.20634920634921),
    @(29, 1.498661560125617, 1.599568439761847, 1896, 0.341991341991342, 30.09523809523809),
    @(30, 1.490594752770106, 1.600093093644183, 1878, 0.3387445887445887, 29.80952380952381),
    @(31, 1.479211240234813, 1.598840421125847, 1861, 0.3356782106782107, 29.53968253968254),
    @(32, 1.47360535950891, 1.598619701150234, 1849, 0.3335137085137085, 29.34920634920635),
    @(33, 1.458620806956018, 1.599161987476047, 1834, 0.3308080808080808, 29.11111111111111),
    @(34, 1.414671191863435, 1.599417319372034, 1826, 0.3293650793650794, 28.98412698412698),
    @(35, 1.402161847727561, 1.597280855746159, 1814, 0.3272005772005772, 28.79365079365079),
    @(36, 1.404133431441432, 1.603954165236823, 1782, 0.3214285714285715, 28.28571428571428),
    @(37, 1.384411121512844, 1.60679823718064, 1763, 0.318001443001443, 27.98412698412698),
    @(38, 1.394247834019988, 1.610331931373487, 1733, 0.3125901875901876, 27.50793650793651),
    @(39, 1.336371914299397, 1.605934063155811, 1711, 0.3086219336219336, 27.15873015873016),
    @(40, 1.341725841084197, 1.603339564800934, 1676, 0.3023088023088023, 26.6031746031746),
    @(41, 1.382697046671422, 1.607803705410245, 1619, 0.292027417027417, 25.6984126984127),
    @(42, 1.365965063689669, 1.603424786843718, 1548, 0.2792207792207792, 24.57142857142857),
    @(43, 1.397086470997752, 1.560339786886919, 1472, 0.2655122655122655, 23.36507936507936),
    @(44, 1.437046457796114, 1.501951814983856, 1392, 0.2510822510822511, 22.09523809523809),
    @(45, 1.376052797621871, 1.428046864211716, 1235, 0.2227633477633478, 19.6031746031746),
    @(46, 1.317397275463974, 1.315879943192129, 1019, 0.1838023088023088, 16.17460317460317),
    @(47, 1.057769277428422, 1.215743157886324, 734, 0.1323953823953824, 11.65079365079365),
    @(48, 0.895241151850083, 1.137423283492765, 451, 0.08134920634920635, 7.158730158730159),
    @(49, 0.7065746747268984, 1.148266318916395, 194, 0.03499278499278499, 3.07936507936508),
    @(50, 0.7394028611841493, 1.401531649329043, 48, 0.008658008658008658, 0.7619047619047619)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}

Write-Output "Applied $($rows.Count) row updates."